# Open the two sheets of the loan-product workbook
$wb = $excel.ActiveWorkbook
$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Correct the product short name: "200MS-..." -> "200-MS-..." (missing dash added)
# The same shared string value is used on both sheets (cell B1)
$wsInput.Range("B1").Value  = "200-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"
$wsOutput.Range("B1").Value = "200-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

# Leave a recorded selection of B1 on the input sheet
$wsInput.Range("B1").Select()

# Switch the active sheet to ProductLoanOutput and select B1 there
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
